$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work Breakdown Structure (WBS): prepend task IDs to the relevant backlog items
$ws.Range("B4").Value = "(ID: 'A') Setup the GanttProject. This includes forking the repository from GitHub"
$ws.Range("B7").Value = "(ID: 'B') Identify pinpoint code smells used to design this tool"
$ws.Range("B10").Value = "(ID: 'C') Identify GoF Design patterns used to design this tool"
$ws.Range("B13").Value = "(ID: 'D') Each team member should review three other colleague´s code smells."
$ws.Range("B16").Value = "(ID: 'E') Each team member should review three other colleague´s design patterns."

# Update Scrum board: move the "Everyone reviews three code smells..." card
# from the TO DO column (F16) to the DOING column (H16), keeping the
# destination cell's merged range (H16:I18) and matching format intact.
$ws.Range("H16").Value = $ws.Range("F16").Value2
$ws.Range("F16").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("F16").Value = $null
